$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header (before the current row 2),
# shifting all existing data down by two rows.
$ws.Rows.Item(2).Resize(2).Insert()
$ws.Range("A2:C3").ClearFormats()

# Populate the two newly inserted rows with the new accelerometer readings.
$ws.Cells.Item(2, 1).Value = -1.97176456451416
$ws.Cells.Item(2, 2).Value = 1.745009422302246
$ws.Cells.Item(2, 3).Value = 0.4838592410087585

$ws.Cells.Item(3, 1).Value = -1.845728397369385
$ws.Cells.Item(3, 2).Value = 1.672563552856445
$ws.Cells.Item(3, 3).Value = 0.5211508870124817
